# update sweep dyn master
$wb = $excel.ActiveWorkbook

# RUN_LIST: nb_step_list for run 2 changes from 20000 to 10000
$wsRunList = $wb.Worksheets.Item("RUN_LIST")
$wsRunList.Range("B3").Value = 10000

# Make RUN_LIST the active/selected sheet (was FORCING_DELTA before),
# with the cursor resting on B4.
$wsRunList.Select()
$wsRunList.Range("B4").Select()
